# Edit script: add "year" column with 2010 data for the existing rows, add a
# full new block of 2020 data for the same 16 states, and tidy up the
# left-alignment formatting that was applied to the original table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New "year" header + 2010 values for the existing table (rows 2-17) ---
$ws.Range("C1").Value = "year"

$years2010 = New-Object 'object[,]' 16,1
for ($i = 0; $i -lt 16; $i++) { $years2010[$i,0] = 2010 }
$ws.Range("C2:C17").Value = $years2010

# --- 2. New block of 2020 data for the same states (rows 18-33) ---
$data2020 = New-Object 'object[,]' 16,2
$data2020[0,0]  = "Johor";            $data2020[0,1]  = 77.4
$data2020[1,0]  = "Kedah";            $data2020[1,1]  = 67.3
$data2020[2,0]  = "Kelantan";         $data2020[2,1]  = 44.1
$data2020[3,0]  = "Melaka";           $data2020[3,1]  = 90.9
$data2020[4,0]  = "Negeri Sembilan";  $data2020[4,1]  = 69.3
$data2020[5,0]  = "Pahang";           $data2020[5,1]  = 52.8
$data2020[6,0]  = "Pulau Pinang";     $data2020[6,1]  = 92.5
$data2020[7,0]  = "Perak";            $data2020[7,1]  = 72
$data2020[8,0]  = "Perlis";           $data2020[8,1]  = 53.8
$data2020[9,0]  = "Selangor";         $data2020[9,1]  = 96.8
$data2020[10,0] = "Terengganu";       $data2020[10,1] = 64.2
$data2020[11,0] = "Sabah";            $data2020[11,1] = 54.7
$data2020[12,0] = "Sarawak";          $data2020[12,1] = 57
$data2020[13,0] = "Kuala Lumpur";     $data2020[13,1] = 100
$data2020[14,0] = "Labuan";           $data2020[14,1] = 88.9
$data2020[15,0] = "Putrajaya";        $data2020[15,1] = 100

$ws.Range("A18:B33").Value = $data2020

$years2020 = New-Object 'object[,]' 16,1
for ($i = 0; $i -lt 16; $i++) { $years2020[$i,0] = 2020 }
$ws.Range("C18:C33").Value = $years2020

# Give the new "thousands"-style cells for 2020 the same accounting number
# format used for the matching rows in the 2010 block, so the two blocks
# look consistent (Johor, Kelantan, N.Sembilan, P.Pinang, Perlis, Terengganu).
$acctFormat = '_-* #,##0.0_-;\-* #,##0.0_-;_-* "-"??_-;_-@'
$ws.Range("B18").NumberFormat = $acctFormat
$ws.Range("B20").NumberFormat = $acctFormat
$ws.Range("B22").NumberFormat = $acctFormat
$ws.Range("B24").NumberFormat = $acctFormat
$ws.Range("B26").NumberFormat = $acctFormat
$ws.Range("B28").NumberFormat = $acctFormat

# --- 3. Remove the explicit "left" horizontal alignment that used to be on
#        every cell of the table; the cleaned-up sheet just uses the
#        default/general alignment everywhere. ---
$ws.Range("A1:C33").HorizontalAlignment = 1

# --- 4. Misc cosmetic bits that came along with the cleanup ---
$ws.Range("F8").Select()
